# Auto-generated edit script
# Updates 'want-to-go count' (F) and 'lowest price' (G) figures across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 164
$ws.Range("F5").Value = 48
$ws.Range("F6").Value = 2740
$ws.Range("F7").Value = 87
$ws.Range("F8").Value = 1622
$ws.Range("G8").Value = 70
$ws.Range("F9").Value = 7420
$ws.Range("F11").Value = 7606
$ws.Range("F13").Value = 33
$ws.Range("F15").Value = 6107
$ws.Range("F16").Value = 3245
$ws.Range("F17").Value = 3618
$ws.Range("F18").Value = 13
$ws.Range("F19").Value = 7
$ws.Range("F20").Value = 13
$ws.Range("F22").Value = 442
$ws.Range("F23").Value = 4
$ws.Range("F24").Value = 280
$ws.Range("F25").Value = 279
$ws.Range("F26").Value = 2109
$ws.Range("F28").Value = 336
$ws.Range("F29").Value = 924
$ws.Range("F31").Value = 1076
$ws.Range("F33").Value = 15
$ws.Range("F34").Value = 2597
$ws.Range("F35").Value = 1447
$ws.Range("F36").Value = 7
$ws.Range("F37").Value = 12
$ws.Range("F38").Value = 19
$ws.Range("F39").Value = 3220
$ws.Range("F40").Value = 150
$ws.Range("F41").Value = 239
$ws.Range("F42").Value = 31
$ws.Range("F43").Value = 895
$ws.Range("F44").Value = 475
$ws.Range("F45").Value = 1259
$ws.Range("F48").Value = 584

$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 55
$ws.Range("F5").Value = 232
$ws.Range("F9").Value = 397
$ws.Range("F10").Value = 29

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 118

$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 164
$ws.Range("F6").Value = 55
$ws.Range("F7").Value = 48
$ws.Range("F8").Value = 118
$ws.Range("F9").Value = 2740
$ws.Range("F10").Value = 1622
$ws.Range("G10").Value = 70
$ws.Range("F11").Value = 232
$ws.Range("F13").Value = 7420
$ws.Range("F14").Value = 7606
$ws.Range("F17").Value = 6107
$ws.Range("F18").Value = 3245
$ws.Range("F19").Value = 3618
$ws.Range("F20").Value = 13
$ws.Range("F21").Value = 7
$ws.Range("F23").Value = 442
$ws.Range("F24").Value = 4
$ws.Range("F25").Value = 29
$ws.Range("F26").Value = 280
$ws.Range("F28").Value = 279
$ws.Range("F29").Value = 2109
$ws.Range("F34").Value = 336
$ws.Range("F35").Value = 924
$ws.Range("F37").Value = 15
$ws.Range("F38").Value = 2597
$ws.Range("F39").Value = 1447
$ws.Range("F40").Value = 7
$ws.Range("F41").Value = 12
$ws.Range("F43").Value = 3220
$ws.Range("F44").Value = 239
$ws.Range("F45").Value = 895
$ws.Range("F46").Value = 475
$ws.Range("F47").Value = 1259
$ws.Range("F49").Value = 584
